# "Weighting & Scaling update & heatmap"
#
# - Scaling!C2:C4 switch from hard-coded 100s to =MAX(ecological_params!..)
#   formulas that pull the real maxima from the ecological_params sheet.
# - Two new helper columns (E "Optimal", F "Threshold") are added next to
#   the existing Scaling table, pre-formatted with the same bordered look
#   as the rest of the table (ready for the heatmap values).
# - The "Scaling" tab becomes the active/selected tab (selection sitting on
#   H12, near the new heatmap area) instead of "ecological_params".

$wb  = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Scaling")

# --- New "Optimal" / "Threshold" header cells (E1:F1), same look as D1 ---
$ws2.Range("D1").Copy()
$ws2.Range("E1:F1").PasteSpecial(-4122)
$ws2.Range("E1").Value = "Optimal"
$ws2.Range("F1").Value = "Threshold"

# --- Column C now derives the scaling max from ecological_params ---
$ws2.Range("C2").Formula = "=MAX(ecological_params!B2:D2)"
$ws2.Range("C3").Formula = "=MAX(ecological_params!B3:D3)"
$ws2.Range("C4").Formula = "=MAX(ecological_params!B4:D4)"

# --- New (currently empty) bordered heatmap cells E2:F4 ---
$ws2.Range("E2:F4").Borders.LineStyle = 1

# --- Make "Scaling" the active sheet / selection, matching the new view ---
$ws2.Activate()
$ws2.Range("H12").Select()
